# Adapt xlspython results sheet to the new MVC architecture:
# drop the "Prénom" / "État" columns (old A:B) and the trailing
# "Réponse 1" column (old E), leaving only "Temps utilisé" and
# "Note/10,00" which become the new A:B columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Prénom" and "État" columns; this shifts
# "Temps utilisé" and "Note/10,00" (with their cell styles) into A:B.
$ws.Columns("A:B").Delete()

# Remove the trailing "Réponse 1" column, now at C.
$ws.Columns("C:C").Delete()
